$d = $word.ActiveDocument

$replacements = @(
    @{old="965×8="; new="717×5="},
    @{old="182×8="; new="210×7="},
    @{old="615×6="; new="368×4="},
    @{old="963×8="; new="691×7="},
    @{old="423×7="; new="429×9="},
    @{old="889×3="; new="701×3="},
    @{old="196×3="; new="310×4="},
    @{old="490×2="; new="243×9="},
    @{old="448×6="; new="298×9="},
    @{old="139×3="; new="408×2="},
    @{old="878×2="; new="617×7="},
    @{old="756×5="; new="399×9="},
    @{old="754×9="; new="713×5="},
    @{old="323×8="; new="398×4="},
    @{old="142×9="; new="194×5="},
    @{old="532×3="; new="836×6="},
    @{old="584×9="; new="113×9="},
    @{old="957×2="; new="133×2="},
    @{old="587×4="; new="955×5="},
    @{old="318×5="; new="880×3="},
    @{old="291×2="; new="110×8="},
    @{old="486×6="; new="779×8="},
    @{old="486×2="; new="487×8="},
    @{old="849×9="; new="361×5="},
    @{old="511×6="; new="638×7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
